# Added IP checksum routine
# --------------------------------------------------------------
# 1. Structural changes: insert a new "offset (hex)" column after
#    the existing offset column, and insert 5 new rows for the
#    "pseudo header for checksum only" block used by the new
#    UDP/IP checksum section.
# --------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new column B (old B "length" -> C, old C "field" -> D)
$ws.Columns("B:B").Insert()

# Insert 4 new rows before the old row 32 (now still row 32, because
# the new column insert does not move rows). Old row 32 needs to land
# on row 36, a shift of +4 -- the previously-blank row 31 shifts down
# to become row 35, which we then also populate as part of the new
# 5-row checksum block (31-35).
$ws.Rows("31:34").Insert()

# --------------------------------------------------------------
# 2. New column header cell for row 1 (keeps title row formatting)
# --------------------------------------------------------------
$ws.Range("B1").Font.Bold = $true

# --------------------------------------------------------------
# 3. Header row (row 3): add "offset (hex)" header between
#    "offset" and "length"
# --------------------------------------------------------------
$ws.Range("B3").Value = "offset (hex)"
$ws.Range("B3").HorizontalAlignment = -4108

# --------------------------------------------------------------
# 4. New "offset (hex)" formula column: DEC2HEX(A,2) for every row
#    that has an offset value in column A. The sheet has several
#    blank separator rows, so this is done in blocks.
# --------------------------------------------------------------
$ws.Range("B4").Formula = "=DEC2HEX(A4,2)"
$ws.Range("B5:B6").Formula = "=DEC2HEX(A5,2)"
$ws.Range("B10:B19").Formula = "=DEC2HEX(A10,2)"
$ws.Range("B23:B28").Formula = "=DEC2HEX(A23,2)"
$ws.Range("B31:B40").Formula = "=DEC2HEX(A31,2)"

# --------------------------------------------------------------
# 5. New "IP checksum" pseudo-header rows (31-35), italic styling.
# --------------------------------------------------------------

# A31:A34 share one relative formula (next offset minus this length)
$ws.Range("A31:A34").Formula = "=A32-C31"
$ws.Range("A35").Formula = "=A36-C35"

$ws.Range("C31").Value = 4
$ws.Range("C32").Value = 4
$ws.Range("C33").Value = 1
$ws.Range("C34").Value = 1
$ws.Range("C35").Value = 2

$ws.Range("D31").Value = "Source IP address"
$ws.Range("D32").Value = "Destination IP address"
$ws.Range("D33").Value = "Zero byte"
$ws.Range("D34").Value = "Protocol"
$ws.Range("D35").Value = "Length (as below)"

$ws.Range("E31").Value = "Pseudo header for checksum only"
$ws.Range("E32").Value = """"
$ws.Range("E33").Value = """"
$ws.Range("E34").Value = """"
$ws.Range("E35").Value = """"

# Italic styling for the whole new block (A31:E35), centred for A:C
$checksumBlock = $ws.Range("A31:E35")
$checksumBlock.Font.Italic = $true

$abc = $ws.Range("A31:C35")
$abc.HorizontalAlignment = -4108

# --------------------------------------------------------------
# 5b. New trailing "Data" row (40), appended after the (shifted)
#     original UDP frame rows, which now end at row 39.
# --------------------------------------------------------------
$ws.Range("A40").Formula = "=A39+C39"
$ws.Range("D40").Value = "Data"

# --------------------------------------------------------------
# 6. Column widths: the new "offset (hex)" column gets its own width.
# --------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 12.43

# --------------------------------------------------------------
# 7. Selection, matching the saved workbook state in the diff.
# --------------------------------------------------------------
$ws.Range("F22").Select()
